$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 119.0815153333333
$ws.Range("H2").Value = 357.244546
$ws.Range("I2").Value = 0.431812569872284
$ws.Range("J2").Value = 0.4318125698722839
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.055739
$ws.Range("N2").Value = 48.167217
$ws.Range("O2").Value = 0.3924791516302356
$ws.Range("P2").Value = 0.3924791516302356
$ws.Range("Q2").Value = 1911.941729916498
$ws.Range("R2").Value = 17207.47556924848
$ws.Range("S2").Value = 0.1694774310867459
$ws.Range("T2").Value = 0.1694774310867458
$ws.Range("G3").Value = 119.0815153333333
$ws.Range("H3").Value = 357.244546
$ws.Range("I3").Value = 0.431812569872284
$ws.Range("J3").Value = 0.4318125698722839
$ws.Range("O3").Value = 0.3954240805436893
$ws.Range("P3").Value = 0.3954240805436893
$ws.Range("Q3").Value = 1926.287797619412
$ws.Range("R3").Value = 17336.59017857471
$ws.Range("S3").Value = 0.1707490884089555
$ws.Range("T3").Value = 0.1707490884089554
$ws.Range("G4").Value = 119.0815153333333
$ws.Range("H4").Value = 357.244546
$ws.Range("I4").Value = 0.431812569872284
$ws.Range("J4").Value = 0.4318125698722839
$ws.Range("M4").Value = 8.676563666666667
$ws.Range("N4").Value = 26.029691
$ws.Range("O4").Value = 0.2120967678260751
$ws.Range("P4").Value = 0.2120967678260751
$ws.Range("Q4").Value = 1033.21834931281
$ws.Range("R4").Value = 9298.965143815287
$ws.Range("S4").Value = 0.09158605037658267
$ws.Range("T4").Value = 0.09158605037658264
$ws.Range("I5").Value = 0.4460879372303943
$ws.Range("J5").Value = 0.4460879372303942
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 16.055739
$ws.Range("N5").Value = 48.167217
$ws.Range("O5").Value = 0.3924791516302356
$ws.Range("P5").Value = 0.3924791516302356
$ws.Range("Q5").Value = 1975.148946348228
$ws.Range("R5").Value = 17776.34051713405
$ws.Range("S5").Value = 0.1750802151566669
$ws.Range("T5").Value = 0.1750802151566669
$ws.Range("I6").Value = 0.4460879372303943
$ws.Range("J6").Value = 0.4460879372303942
$ws.Range("O6").Value = 0.3954240805436893
$ws.Range("P6").Value = 0.3954240805436893
$ws.Range("S6").Value = 0.1763939124209597
$ws.Range("T6").Value = 0.1763939124209596
$ws.Range("I7").Value = 0.4460879372303943
$ws.Range("J7").Value = 0.4460879372303942
$ws.Range("M7").Value = 8.676563666666667
$ws.Range("N7").Value = 26.029691
$ws.Range("O7").Value = 0.2120967678260751
$ws.Range("P7").Value = 0.2120967678260751
$ws.Range("Q7").Value = 1067.375695640044
$ws.Range("R7").Value = 9606.381260760396
$ws.Range("S7").Value = 0.09461380965276772
$ws.Range("T7").Value = 0.09461380965276769
$ws.Range("G8").Value = 33.50679633333333
$ws.Range("H8").Value = 100.520389
$ws.Range("I8").Value = 0.1215021138451521
$ws.Range("J8").Value = 0.121502113845152
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 16.055739
$ws.Range("N8").Value = 48.167217
$ws.Range("O8").Value = 0.3924791516302356
$ws.Range("P8").Value = 0.3924791516302356
$ws.Range("Q8").Value = 537.976376654157
$ws.Range("R8").Value = 4841.787389887413
$ws.Range("S8").Value = 0.04768704656322559
$ws.Range("T8").Value = 0.04768704656322557
$ws.Range("G9").Value = 33.50679633333333
$ws.Range("H9").Value = 100.520389
$ws.Range("I9").Value = 0.1215021138451521
$ws.Range("J9").Value = 0.121502113845152
$ws.Range("O9").Value = 0.3954240805436893
$ws.Range("P9").Value = 0.3954240805436893
$ws.Range("Q9").Value = 542.0130297598906
$ws.Range("R9").Value = 4878.117267839015
$ws.Range("S9").Value = 0.04804486165133392
$ws.Range("T9").Value = 0.0480448616513339
$ws.Range("G10").Value = 33.50679633333333
$ws.Range("H10").Value = 100.520389
$ws.Range("I10").Value = 0.1215021138451521
$ws.Range("J10").Value = 0.121502113845152
$ws.Range("M10").Value = 8.676563666666667
$ws.Range("N10").Value = 26.029691
$ws.Range("O10").Value = 0.2120967678260751
$ws.Range("P10").Value = 0.2120967678260751
$ws.Range("Q10").Value = 290.7238516521999
$ws.Range("R10").Value = 2616.514664869799
$ws.Range("S10").Value = 0.02577020563059257
$ws.Range("T10").Value = 0.02577020563059256
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.16474
$ws.Range("H11").Value = 0.49422
$ws.Range("I11").Value = 0.000597379052169715
$ws.Range("J11").Value = 0.000597379052169715
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 16.055739
$ws.Range("N11").Value = 48.167217
$ws.Range("O11").Value = 0.3924791516302356
$ws.Range("P11").Value = 0.3924791516302356
$ws.Range("Q11").Value = 2.64502244286
$ws.Range("R11").Value = 23.80520198574
$ws.Range("S11").Value = 0.000234458823597244
$ws.Range("T11").Value = 0.000234458823597244
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.16474
$ws.Range("H12").Value = 0.49422
$ws.Range("I12").Value = 0.000597379052169715
$ws.Range("J12").Value = 0.000597379052169715
$ws.Range("O12").Value = 0.3954240805436893
$ws.Range("P12").Value = 0.3954240805436893
$ws.Range("Q12").Value = 2.664869109966667
$ws.Range("R12").Value = 23.9838219897
$ws.Range("S12").Value = 0.0002362180624402702
$ws.Range("T12").Value = 0.0002362180624402702
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.16474
$ws.Range("H13").Value = 0.49422
$ws.Range("I13").Value = 0.000597379052169715
$ws.Range("J13").Value = 0.000597379052169715
$ws.Range("M13").Value = 8.676563666666667
$ws.Range("N13").Value = 26.029691
$ws.Range("O13").Value = 0.2120967678260751
$ws.Range("P13").Value = 0.2120967678260751
$ws.Range("Q13").Value = 1.429377098446667
$ws.Range("R13").Value = 12.86439388602
$ws.Range("S13").Value = 0.0001267021661322009
$ws.Range("T13").Value = 0.0001267021661322009
